# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet (positioned right after "总计" and
# before "2021-Q2"), populates it with the quarterly fund-holding table,
# and adds a corresponding summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet before the existing "2021-Q2" sheet.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet    = $wb.Worksheets.Item(2)
$newSheet   = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2. Populate the header row of the new sheet (bold / bordered / centered
#    style, matching the other quarterly sheets).
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $newSheet.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
}
$newSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Fill in the fund data rows (2-9). Columns B-G hold text values
#    (fund codes / formatted numbers kept as strings); column H is numeric.
# ---------------------------------------------------------------------
$q4Data = @(
    @("000179", "广发美国房地产指数（QDII）人民币A", "1.60", "92.49", "10.07", "0.1611", 1),
    @("000180", "广发美国房地产指数（QDII）美元A", "1.60", "92.49", "10.07", "0.1611", 1),
    @("160140", "南方道琼斯美国精选REIT指数（QDII-LOF）A", "0.80", "92.31", "10.66", "0.0853", 1),
    @("160141", "南方道琼斯美国精选REIT指数（QDII-LOF）C", "0.39", "92.31", "10.66", "0.0416", 1),
    @("070031", "嘉实全球房地产（QDII）", "0.39", "94.72", "8.06", "0.0314", 1),
    @("320017", "诺安全球收益不动产（QDII）", "0.24", "68.42", "8.64", "0.0207", 1),
    @("016278", "广发美国房地产指数（QDII）人民币C", "0.01", "92.49", "10.07", "0.0010", 1),
    @("016279", "广发美国房地产指数（QDII）美元C", "0.01", "92.49", "10.07", "0.0010", 1)
)

$startRow = 2
$newSheet.Range("A2:A9").NumberFormat = "General"
$newSheet.Range("B2:G9").NumberFormat = "@"

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $q4Data[$i]

    $newSheet.Cells.Item($r, 1).Value = $i          # A: 0-based index
    $newSheet.Cells.Item($r, 2).Value = $rowVals[0] # B: fund code
    $newSheet.Cells.Item($r, 3).Value = $rowVals[1] # C: fund name
    $newSheet.Cells.Item($r, 4).Value = $rowVals[2] # D: fund size
    $newSheet.Cells.Item($r, 5).Value = $rowVals[3] # E: total stock position
    $newSheet.Cells.Item($r, 6).Value = $rowVals[4] # F: position ratio
    $newSheet.Cells.Item($r, 7).Value = $rowVals[5] # G: held market value
    $newSheet.Cells.Item($r, 8).Value = $rowVals[6] # H: position rank
}

# Match the "A" column style (bold/centered/bordered) used on the other
# quarterly sheets.
$totalSheet.Range("A2").Copy()
$newSheet.Range("A2:A9").PasteSpecial(-4122)
for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $r = $startRow + $i
    $newSheet.Cells.Item($r, 1).Value = $i
}

$newSheet.PageSetup.LeftMargin   = 0.75 * 72
$newSheet.PageSetup.RightMargin  = 0.75 * 72
$newSheet.PageSetup.TopMargin    = 1 * 72
$newSheet.PageSetup.BottomMargin = 1 * 72
$newSheet.PageSetup.HeaderMargin = 0.5 * 72
$newSheet.PageSetup.FooterMargin = 0.5 * 72

# ---------------------------------------------------------------------
# 4. Insert the 2022-Q4 summary row into "总计" (shifting 2021-Q2 /
#    2021-Q1 down by one row) and renumber the index column.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 8
$totalSheet.Range("D2").Value = 0.5

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q2"
$totalSheet.Range("C3").Value = 9
$totalSheet.Range("D3").Value = 1.96

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q1"
$totalSheet.Range("C4").Value = 12
$totalSheet.Range("D4").Value = 1.41

# Row-insert copies row 1's (header) formatting onto the new row 2 for
# columns B:D - clear that back to the unformatted state the other data
# rows use, keeping only column A's existing numbered style.
$totalSheet.Range("B2:D2").ClearFormats()

# ---------------------------------------------------------------------
# 5. Restore the original active-sheet selection ("总计").
# ---------------------------------------------------------------------
$totalSheet.Activate()
